## Add an "is_multiple_choice" column to the trivia questions sheet.
## The new column is inserted before the existing "answers" column (old F),
## which (together with the old "answer_message" column, old G) shifts one
## slot to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; everything currently in F/G shifts to G/H.
$ws.Columns.Item(6).Insert()

# Header for the new column.
$ws.Cells.Item(1, 6).Value = "is_multiple_choice"

# All data rows (2-10) are marked as multiple choice questions.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Value = "y"
}

# Match the column widths from the target layout: the new column is
# narrower than the two columns that got pushed over.
$ws.Columns.Item(6).ColumnWidth = 18.33
$ws.Columns.Item(9).ColumnWidth = 17.42

# Update the view to match: scrolled right a bit, selection on F7.
$ws.Range("F7").Select()
